$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/6/2023  Through  2/12/2023"

# --- Cells changing from numeric to the shared text "0" (style 14) ---
# Copy a stable same-valued text "0" cell (row 14, untouched by this revision)
# onto the target so both value/type AND style index land correctly.
$ws.Range("C14").Copy($ws.Range("F15"))
$ws.Range("C14").Copy($ws.Range("C18"))

# --- Cells changing from the shared text "0" to a real numeric value (style 16) ---
$ws.Range("I15").Copy($ws.Range("C26"))

# --- Plain numeric value updates ---
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -46.666666666666
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 23
$ws.Range("K16").Value = -47.826086956521
$ws.Range("M16").Value = -71.428571428571
$ws.Range("N16").Value = -92.5
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 16.666666666666
$ws.Range("I17").Value = 39
$ws.Range("J17").Value = 22
$ws.Range("K17").Value = 77.272727272727
$ws.Range("L17").Value = 21.875
$ws.Range("M17").Value = 62.5
$ws.Range("N17").Value = 18.181818181818
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -50
$ws.Range("J18").Value = 13
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 8.333333333333
$ws.Range("M18").Value = -53.571428571428
$ws.Range("N18").Value = -93.121693121693
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 200
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 3.846153846153
$ws.Range("I19").Value = 43
$ws.Range("J19").Value = 52
$ws.Range("K19").Value = -17.307692307692
$ws.Range("L19").Value = 48.275862068965
$ws.Range("M19").Value = 34.375
$ws.Range("N19").Value = -40.277777777777
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = -14.285714285714
$ws.Range("I20").Value = 34
$ws.Range("J20").Value = 35
$ws.Range("K20").Value = -2.857142857142
$ws.Range("L20").Value = 70
$ws.Range("M20").Value = 3.030303030303
$ws.Range("N20").Value = -93.486590038314
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 17.647058823529
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 86
$ws.Range("H21").Value = -10.465116279069
$ws.Range("I21").Value = 142
$ws.Range("J21").Value = 146
$ws.Range("K21").Value = -2.739726027397
$ws.Range("L21").Value = 35.238095238095
$ws.Range("M21").Value = -11.25
$ws.Range("N21").Value = -85.495403472931
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -7.894736842105
$ws.Range("F24").Value = 115
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = 6.481481481481
$ws.Range("I24").Value = 160
$ws.Range("J24").Value = 154
$ws.Range("K24").Value = 3.896103896103
$ws.Range("L24").Value = 10.344827586206
$ws.Range("M24").Value = 53.846153846153
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -6.666666666666
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = 12.5
$ws.Range("I25").Value = 64
$ws.Range("J25").Value = 55
$ws.Range("K25").Value = 16.363636363636
$ws.Range("L25").Value = 72.972972972973
$ws.Range("M25").Value = -9.859154929577
$ws.Range("I26").Value = 2
$ws.Range("K26").Value = -33.333333333333
$ws.Range("L26").Value = -50
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 7
$ws.Range("H27").Value = 40
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 28.571428571428
